$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.998.90"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'0.7918"
$ws.Range("E5").Value = "  +4.20%  "
$ws.Range("D6").Value = "'241.86"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D8").Value = "'0.3154"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").Value = "'26.25"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "'0.06905"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "'0.08000"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.902.51"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'0.7417"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'5.188"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "'92.99"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "29.995.21"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'5.861"
$ws.Range("E18").Value = "  -5.67%  "
$ws.Range("D19").Value = "'245.58"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("D20").Value = "'0.000007738"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "2.151.70"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'6.831"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").Value = "'168.03"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'9.231"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("D28").Value = "'18.91"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'1.514"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'4.310"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'0.05516"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "'1.257"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "'0.7319"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "'2.723"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'0.01925"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "'2.782"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.4411"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.103"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "'72.21"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'0.8372"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'1.872"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'100.50"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.529"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "'988.26"
$ws.Range("E48").Value = "  +6.81%  "
$ws.Range("D49").Value = "2.058.36"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'36.21"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'0.05936"
$ws.Range("E51").Value = "  -0.56%  "

Write-Host "Applied 98 cell updates"